$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Delete the data row for account 005646524 / EVANGELINA (row 2, right below the
# header row). This removes the whole row and shifts all subsequent rows up by one.
$ws.Rows.Item(2).Delete()
